# "Generate Report for Archive"
#
# This CI-generated localization-status report is re-emitted with its
# Status vocabulary gaining a new recognized value, "In Translation",
# that is registered in the workbook's string table (it is not yet used
# by any row in this particular snapshot -- none of the existing
# File Name / Path / Status / Date values for the Overview, zh-cn or
# de-de sheets change). Record the new status value on the Overview
# sheet, just below the existing summary table, so the archive keeps a
# record of every known status without disturbing any of the existing
# report data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A7").Value = "In Translation"
